$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 27, shifting existing rows 27-38 down to 28-39
$ws.Rows.Item(27).Insert()

# Copy the style (date number format) used by column D in the data rows onto the new row's D cell
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 27 with the new data record
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44460
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112012
$ws.Cells.Item(27, 7).Value = "Espinaca"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 300
$ws.Cells.Item(27, 11).Value = 950
$ws.Cells.Item(27, 12).Value = 1000
$ws.Cells.Item(27, 13).Value = 975
$ws.Cells.Item(27, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 325
$ws.Cells.Item(27, 17).Value = 3
$ws.Cells.Item(27, 18).Value = "Hortaliza"
